# Add an "Expiration " column (G) to the food data sheet, with a default
# expiration value of 10 for every existing food row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header
$ws.Range("G1").Value = "Expiration "

# Fill in the expiration value for each of the 9 data rows (rows 2-10)
for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 7).Value = 10
}

# Resize the new/adjacent columns to fit their content, like Excel does
# automatically when you widen a column after typing into it.
$ws.Columns("F:G").AutoFit()

# Leave the selection on the first data cell of the new column, matching
# where the author's cursor ended up after adding the data.
$ws.Range("G2").Select()
